# DaySale report update: 4 new medicine rows inserted into the alphabetically
# sorted item list, the running total recalculated, and the generated-on
# timestamp bumped by one minute.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Capture the data that currently lives in the rows we are about to shift
#    down (the tail of the item list, the totals row and the footer row),
#    before anything gets overwritten.
# ---------------------------------------------------------------------------

# Old rows 14..17 -> new rows 18..21 (WATER FOR INJECTION .. معجون سيجنال)
$tailRows = @()
for ($r = 14; $r -le 17; $r++) {
    $rowVals = @{}
    for ($col = 1; $col -le 17; $col++) {
        $rowVals[$col] = $ws.Cells.Item($r, $col).Value2
    }
    $tailRows += ,$rowVals
}

# Old totals row (18) -> new row 22
$totalValue = $ws.Cells.Item(18, 16).Value2

# Old footer row (19) -> new row 23
$footerVals = @{}
for ($col = 1; $col -le 17; $col++) {
    $footerVals[$col] = $ws.Cells.Item(19, $col).Value2
}

# ---------------------------------------------------------------------------
# 2. Extend the sheet's item-row formatting down to the four freshly needed
#    rows (18-21) by copying the format of an existing item row. This avoids
#    Rows.Insert(), which mints brand-new (border-less) style records.
# ---------------------------------------------------------------------------

$ws.Range("A17:Q17").Copy() | Out-Null
$ws.Range("A18:Q21").PasteSpecial(-4122) | Out-Null

# Re-home the totals row format (old row 18) onto new row 22, and the footer
# row format (old row 19) onto new row 23.
$ws.Range("A18:Q18").Copy() | Out-Null
$ws.Range("A22:Q22").PasteSpecial(-4122) | Out-Null

$ws.Range("A19:Q19").Copy() | Out-Null
$ws.Range("A23:Q23").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Write the new item list (7..21), interleaving the 4 new medicines at
#    their alphabetically-correct spots with the shifted-down existing rows.
# ---------------------------------------------------------------------------

function Set-ItemRow($row, $idx, $name, $stock, $limit, $price, $sell, $txn) {
    $ws.Cells.Item($row, 1).Value2 = $idx      # A - م
    $ws.Cells.Item($row, 3).Value2 = $name     # C - الاسم
    $ws.Cells.Item($row, 8).Value2 = $stock    # H - الرصيد الحالي
    $ws.Cells.Item($row, 12).Value2 = $limit   # L - حد الطلب
    $ws.Cells.Item($row, 14).Value2 = $price   # N - السعر
    $ws.Cells.Item($row, 16).Value2 = $sell    # P - سعر  البيع
    $ws.Cells.Item($row, 17).Value2 = $txn     # Q - عدد التعااملات
}

Set-ItemRow 7  1  "ALFACALCIDOL 2MCG/ML ORAL DPS. 10 ML" "1:0" 1 86   86.0000 "1:0"
Set-ItemRow 8  2  "BI-PROFENID 150MG 20 SCORED TABS." "2:0" 1 54 27.0000 "0:1"
Set-ItemRow 9  3  "CLAVIMOX 642.9MG/5ML PD. FOR ORAL SUSP. 70ML" "1:0" 1 89 89.0000 "1:0"
Set-ItemRow 10 4  "CONTAFEVER N 200MG/5ML SUSP. 120ML" "11:0" 1 33 33.0000 "1:0"
Set-ItemRow 11 5  "COPAD 10.000 30 CAPS." "0:0" 1 215 70.9500 "0:1"
Set-ItemRow 12 6  "CURAM 1GM 12 F.C. TABS." "0:2" 1 182 60.0600 "0:1"
Set-ItemRow 13 7  "DANSET 8MG/4ML 3 AMP." "1:1" 1 142.5 94.0500 "0:2"
Set-ItemRow 14 8  "DIGESTOZYME 20 ENTERIC COATED TABLETS (DIGENORM)" "2:0" 1 62 31.0000 "0:1"
Set-ItemRow 15 9  "EPICOPRED 5 MG 30 ORODISPERSIBLE TABS." "0:2" 1 69 22.7700 "0:1"
Set-ItemRow 16 10 "ERASTAPEX TRIO 5/40/12.5MG 30 F.C. TABS" "1:0" 1 144 432.0000 "3:0"
Set-ItemRow 17 11 "OSTEOCARE 30 TABS" "1:0" 1 150 75.0000 "0:1"
Set-ItemRow 18 12 $tailRows[0][3] $tailRows[0][8] $tailRows[0][12] $tailRows[0][14] $tailRows[0][16] $tailRows[0][17]
Set-ItemRow 19 13 $tailRows[1][3] $tailRows[1][8] $tailRows[1][12] $tailRows[1][14] $tailRows[1][16] $tailRows[1][17]
Set-ItemRow 20 14 $tailRows[2][3] $tailRows[2][8] $tailRows[2][12] $tailRows[2][14] $tailRows[2][16] $tailRows[2][17]
Set-ItemRow 21 15 $tailRows[3][3] $tailRows[3][8] $tailRows[3][12] $tailRows[3][14] $tailRows[3][16] $tailRows[3][17]

# ---------------------------------------------------------------------------
# 4. Totals row (new row 22) and footer row (new row 23).
# ---------------------------------------------------------------------------

$newTotal = [math]::Round($totalValue + 86 + 70.95 + 31 + 75, 2)
$ws.Cells.Item(22, 16).Value2 = $newTotal

$ws.Cells.Item(23, 1).Value2 = "Friday, 5 September, 2025 11:58 AM"
$ws.Cells.Item(23, 7).Value2 = $footerVals[7]
$ws.Cells.Item(23, 11).Value2 = $footerVals[11]

# ---------------------------------------------------------------------------
# 5. Row heights - mirror the target (each row keeps the height belonging to
#    its position, the new rows pick up the repeating 25.5/24.75 pattern).
# ---------------------------------------------------------------------------

$ws.Rows.Item(18).RowHeight = 24.75
$ws.Rows.Item(19).RowHeight = 25.5
$ws.Rows.Item(20).RowHeight = 24.75
$ws.Rows.Item(21).RowHeight = 25.5
$ws.Rows.Item(22).RowHeight = 25.5
$ws.Rows.Item(23).RowHeight = 16.5

# ---------------------------------------------------------------------------
# 6. Merged cells for the newly-populated rows.
# ---------------------------------------------------------------------------

$ws.Range("A18:B18").Merge() | Out-Null
$ws.Range("C18:G18").Merge() | Out-Null
$ws.Range("H18:K18").Merge() | Out-Null
$ws.Range("L18:M18").Merge() | Out-Null
$ws.Range("N18:O18").Merge() | Out-Null

$ws.Range("A19:B19").Merge() | Out-Null
$ws.Range("C19:G19").Merge() | Out-Null
$ws.Range("H19:K19").Merge() | Out-Null
$ws.Range("L19:M19").Merge() | Out-Null
$ws.Range("N19:O19").Merge() | Out-Null

$ws.Range("A20:B20").Merge() | Out-Null
$ws.Range("C20:G20").Merge() | Out-Null
$ws.Range("H20:K20").Merge() | Out-Null
$ws.Range("L20:M20").Merge() | Out-Null
$ws.Range("N20:O20").Merge() | Out-Null

$ws.Range("A21:B21").Merge() | Out-Null
$ws.Range("C21:G21").Merge() | Out-Null
$ws.Range("H21:K21").Merge() | Out-Null
$ws.Range("L21:M21").Merge() | Out-Null
$ws.Range("N21:O21").Merge() | Out-Null

$ws.Range("P22:Q22").Merge() | Out-Null

$ws.Range("A23:F23").Merge() | Out-Null
$ws.Range("G23:I23").Merge() | Out-Null
$ws.Range("K23:Q23").Merge() | Out-Null

Write-Output "DaySale report updated: 4 rows inserted, total and timestamp refreshed."
